$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1026.9166
$ws.Range("I6").Value = 192.4
$ws.Range("K6").Value = 577.2
$ws.Range("M6").Value = -465.2
$ws.Range("H8").Value = 145.8
$ws.Range("I8").Value = 50.88889
$ws.Range("K8").Value = 152.66667
$ws.Range("M8").Value = -13.66667000000001
$ws.Range("H31").Value = 96.666664
$ws.Range("I31").Value = 96.666664
$ws.Range("K31").Value = 289.999992
$ws.Range("M31").Value = -59.99999200000002
$ws.Range("H38").Value = 378.5
$ws.Range("I38").Value = 55
$ws.Range("J38").Value = 1133.3334
$ws.Range("K38").Value = 165
$ws.Range("L38").Value = 3400.0002
$ws.Range("M38").Value = 207
$ws.Range("N38").Value = -4144.0002
$ws.Range("H39").Value = 542.3
$ws.Range("I39").Value = 427.875
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 1283.625
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = -987.625
$ws.Range("N39").Value = -3592
$ws.Range("H43").Value = 3519
$ws.Range("I43").Value = 4500.5
$ws.Range("K43").Value = 4500.5
$ws.Range("M43").Value = -4431.5
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("H86").Value = 5508.4287
$ws.Range("I86").Value = 2999
$ws.Range("K86").Value = 2999
$ws.Range("M86").Value = -1876
$ws.Range("H88").Value = 3284.3572
$ws.Range("J88").Value = 3228
$ws.Range("L88").Value = 3228
$ws.Range("N88").Value = -4040
$ws.Range("H89").Value = 5508.4287
$ws.Range("I89").Value = 2999
$ws.Range("K89").Value = 14995
$ws.Range("M89").Value = -9379
$ws.Range("H91").Value = 3284.3572
$ws.Range("J91").Value = 3228
$ws.Range("L91").Value = 3228
$ws.Range("N91").Value = -6036
$ws.Range("H132").Value = 3097.4285
$ws.Range("I132").Value = 1947.1666
$ws.Range("K132").Value = 5841.4998
$ws.Range("M132").Value = -3311.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H63").Value = 5111
$ws.Range("I63").Value = 4166.5
$ws.Range("K63").Value = 4166.5
$ws.Range("M63").Value = -3480.5
$ws.Range("H66").Value = 5111
$ws.Range("I66").Value = 4166.5
$ws.Range("K66").Value = 20832.5
$ws.Range("M66").Value = -17400.5
$ws.Range("H88").Value = 2034
$ws.Range("I88").Value = 1540.8572
$ws.Range("J88").Value = 2527.1428
$ws.Range("K88").Value = 1540.8572
$ws.Range("L88").Value = 2527.1428
$ws.Range("M88").Value = -1134.8572
$ws.Range("N88").Value = -3339.1428
$ws.Range("H91").Value = 2034
$ws.Range("I91").Value = 1540.8572
$ws.Range("J91").Value = 2527.1428
$ws.Range("K91").Value = 1540.8572
$ws.Range("L91").Value = 2527.1428
$ws.Range("M91").Value = -136.8571999999999
$ws.Range("N91").Value = -5335.1428
$ws.Range("H132").Value = 3221.889
$ws.Range("I132").Value = 3124.625
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 9373.875
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6843.875
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1224.75
$ws.Range("I20").Value = 1299.5
$ws.Range("J20").Value = 1150
$ws.Range("K20").Value = 1299.5
$ws.Range("L20").Value = 1150
$ws.Range("M20").Value = -1052.5
$ws.Range("N20").Value = -1644
$ws.Range("H22").Value = 395
$ws.Range("I22").Value = 458.9091
$ws.Range("J22").Value = 160.66667
$ws.Range("K22").Value = 458.9091
$ws.Range("L22").Value = 160.66667
$ws.Range("M22").Value = -285.9091
$ws.Range("N22").Value = -506.66667
$ws.Range("H26").Value = 45140
$ws.Range("I26").Value = 45140
$ws.Range("K26").Value = 45140
$ws.Range("M26").Value = -44848
$ws.Range("H86").Value = 1448.091
$ws.Range("I86").Value = 1967.7142
$ws.Range("J86").Value = 538.75
$ws.Range("K86").Value = 1967.7142
$ws.Range("L86").Value = 538.75
$ws.Range("M86").Value = -844.7141999999999
$ws.Range("N86").Value = -2784.75
$ws.Range("H89").Value = 1448.091
$ws.Range("I89").Value = 1967.7142
$ws.Range("J89").Value = 538.75
$ws.Range("K89").Value = 9838.571
$ws.Range("L89").Value = 2693.75
$ws.Range("M89").Value = -4222.571
$ws.Range("N89").Value = -13925.75
$ws.Range("H94").Value = 1639.9286
$ws.Range("I94").Value = 1788.25
$ws.Range("K94").Value = 1788.25
$ws.Range("M94").Value = -1337.25
$ws.Range("H96").Value = 11749
$ws.Range("I96").Value = 11749
$ws.Range("K96").Value = 11749
$ws.Range("M96").Value = -9003
$ws.Range("H105").Value = 2038.2
$ws.Range("I105").Value = 1819
$ws.Range("K105").Value = 1819
$ws.Range("M105").Value = -72
$ws.Range("H107").Value = 660.1786
$ws.Range("I107").Value = 658.7037
$ws.Range("K107").Value = 658.7037
$ws.Range("M107").Value = 1261.2963
$ws.Range("H134").Value = 3111.75
$ws.Range("I134").Value = 3111.75
$ws.Range("K134").Value = 9335.25
$ws.Range("M134").Value = -6800.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3186.4
$ws.Range("I58").Value = 3186.4
$ws.Range("K58").Value = 3186.4
$ws.Range("M58").Value = -2983.4
$ws.Range("H132").Value = 2816.2
$ws.Range("I132").Value = 2816.2
$ws.Range("K132").Value = 8448.599999999999
$ws.Range("M132").Value = -5918.599999999999
$ws.Range("H134").Value = 1993.0625
$ws.Range("I134").Value = 1993.0625
$ws.Range("K134").Value = 5979.1875
$ws.Range("M134").Value = -3444.1875
$ws.Range("H136").Value = 3186.4
$ws.Range("I136").Value = 3186.4
$ws.Range("K136").Value = 9559.200000000001
$ws.Range("M136").Value = -7009.200000000001
$ws.Range("H141").Value = 143549.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 143549.25
$ws.Range("K141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("M141").Value = 143549.25
$ws.Range("N141").Value = -153909.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 369.66666
$ws.Range("J41").Value = 750
$ws.Range("L41").Value = 2250
$ws.Range("N41").Value = -2926

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2684.2856
$ws.Range("I80").Value = 2666.6667
$ws.Range("J80").Value = 2697.5
$ws.Range("K80").Value = 2666.6667
$ws.Range("L80").Value = 2697.5
$ws.Range("M80").Value = -1668.6667
$ws.Range("N80").Value = -4693.5
$ws.Range("H83").Value = 2684.2856
$ws.Range("I83").Value = 2666.6667
$ws.Range("J83").Value = 2697.5
$ws.Range("K83").Value = 13333.3335
$ws.Range("L83").Value = 13487.5
$ws.Range("M83").Value = -8341.333500000001
$ws.Range("N83").Value = -23471.5
$ws.Range("H97").Value = 1788.5
$ws.Range("I97").Value = 1501.5
$ws.Range("J97").Value = 2649.5
$ws.Range("K97").Value = 1501.5
$ws.Range("L97").Value = 2649.5
$ws.Range("M97").Value = -1005.5
$ws.Range("N97").Value = -3641.5
$ws.Range("H132").Value = 2209.1428
$ws.Range("I132").Value = 2159.2
$ws.Range("K132").Value = 6477.599999999999
$ws.Range("M132").Value = -3947.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 33999
$ws.Range("J96").Value = 33999
$ws.Range("L96").Value = 33999
$ws.Range("N96").Value = -39491
$ws.Range("H131").Value = 58333.332
$ws.Range("J131").Value = 58333.332
$ws.Range("L131").Value = 58333.332
$ws.Range("N131").Value = -68413.33199999999
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
